# Update loading_percent values for case with 380 kV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 7.77487549862099),
    @("D2", 9.839353882531107),
    @("E2", 12.34823502320366),
    @("F2", 55.07432857900692),
    @("G2", 3.595516286326692),
    @("M2", 28.62232417437256),
    @("N2", 17.09175230921114),
    @("B3", 7.701677463735026),
    @("D3", 9.488196865936814),
    @("E3", 11.73843576557568),
    @("F3", 52.69921568052324),
    @("G3", 3.608381685423935),
    @("M3", 27.33547170258859),
    @("N3", 17.14685100458605),
    @("B4", 7.658230693850197),
    @("D4", 9.272720861691635),
    @("E4", 11.34729184380394),
    @("F4", 51.2217775676208),
    @("G4", 3.616605203599129),
    @("M4", 26.52672657367357),
    @("N4", 17.18323382238964),
    @("B5", 7.640921800496317),
    @("D5", 9.18508494781887),
    @("E5", 11.18378594592479),
    @("F5", 50.61586490879911),
    @("G5", 3.620038933322151),
    @("M5", 26.19295672432217),
    @("N5", 17.1987063318246),
    @("B6", 7.638072147792325),
    @("D6", 9.17054768201948),
    @("E6", 11.1563907231572),
    @("F6", 50.51505149677724),
    @("G6", 3.620614120129297),
    @("M6", 26.13729614252336),
    @("N6", 17.2013146903331),
    @("B7", 7.657995631282563),
    @("D7", 9.271538086524735),
    @("E7", 11.34510325511123),
    @("F7", 51.21362023753313),
    @("G7", 3.616651176309635),
    @("M7", 26.52224155903817),
    @("N7", 17.18343986764043),
    @("B8", 7.749338032401296),
    @("D8", 9.718326804054264),
    @("E8", 12.1414978914945),
    @("F8", 54.25989906503887),
    @("G8", 3.59988566639213),
    @("M8", 28.18274925949739),
    @("N8", 17.11022339960243),
    @("B9", 7.939448642363345),
    @("D9", 10.589993901559),
    @("E9", 13.56733956090898),
    @("F9", 60.04532556653706),
    @("G9", 3.56952973750546),
    @("M9", 31.27338326726317),
    @("N9", 16.98669033276574),
    @("B10", 8.084594525645604),
    @("D10", 11.2210423186294),
    @("E10", 14.52895616044334),
    @("F10", 64.13958119940624),
    @("G10", 3.54869134786648),
    @("M10", 33.42325507962206),
    @("N10", 16.90786261383011),
    @("B11", 8.151546189942804),
    @("D11", 11.50493249365548),
    @("E11", 14.94737620586293),
    @("F11", 65.96151990606293),
    @("G11", 3.539512715240498),
    @("M11", 34.37193269238229),
    @("N11", 16.87452854559766),
    @("B12", 8.177009419745225),
    @("D12", 11.61189767405657),
    @("E12", 15.10306744981748),
    @("F12", 66.64517382580715),
    @("G12", 3.536078854745343),
    @("M12", 34.72675898337721),
    @("N12", 16.86226344462876),
    @("B13", 8.171520923059301),
    @("D13", 11.58888588931474),
    @("E13", 15.06965931162418),
    @("F13", 66.49822192537692),
    @("G13", 3.536816557591994),
    @("M13", 34.65054004322248),
    @("N13", 16.86488912753074),
    @("B14", 8.153638989000322),
    @("D14", 11.51374378330739),
    @("E14", 14.96024031822467),
    @("F14", 66.01789250909809),
    @("G14", 3.539229378056242),
    @("M14", 34.40121400848312),
    @("N14", 16.87351234887657),
    @("B15", 8.142699450162082),
    @("D15", 11.46764483356278),
    @("E15", 14.89285877752124),
    @("F15", 65.72284805806673),
    @("G15", 3.540712713534058),
    @("M15", 34.24791404075589),
    @("N15", 16.87884075451013),
    @("B16", 8.080235721100776),
    @("D16", 11.20241829066745),
    @("E16", 14.50122670519026),
    @("F16", 64.01965726805092),
    @("G16", 3.549297131006578),
    @("M16", 33.36064869877707),
    @("N16", 16.91009143426375),
    @("B17", 8.042136681941956),
    @("D17", 11.03883511876128),
    @("E17", 14.25608341833911),
    @("F17", 62.9640741922843),
    @("G17", 3.554639490824354),
    @("M17", 32.80867556185945),
    @("N17", 16.92990571954162),
    @("B18", 8.020311463081757),
    @("D18", 10.94445142944516),
    @("E18", 14.11329335139593),
    @("F18", 62.35313382901779),
    @("G18", 3.557740678618298),
    @("M18", 32.48844571713675),
    @("N18", 16.94154037120094),
    @("B19", 8.012937687491393),
    @("D19", 10.91244674390194),
    @("E19", 14.0646403921859),
    @("F19", 62.14564271062692),
    @("G19", 3.55879560675992),
    @("M19", 32.37955603054407),
    @("N19", 16.94552069972541),
    @("B20", 8.046183412191828),
    @("D20", 11.05627998183288),
    @("E20", 14.28236478499828),
    @("F20", 63.07683925737064),
    @("G20", 3.55406785758359),
    @("M20", 32.86772020722863),
    @("N20", 16.92777186097015),
    @("B21", 8.158888539816918),
    @("D21", 11.5358300234079),
    @("E21", 14.99245421772055),
    @("F21", 66.15915029621051),
    @("G21", 3.538519548618552),
    @("M21", 34.47456840363066),
    @("N21", 16.87096983277958),
    @("B22", 8.233179068490955),
    @("D22", 11.84607967096215),
    @("E22", 15.44047584969961),
    @("F22", 68.13689892703945),
    @("G22", 3.528601304692728),
    @("M22", 35.49892334027498),
    @("N22", 16.83592927647574),
    @("B23", 8.193478552680597),
    @("D23", 11.68080669359632),
    @("E23", 15.20283210417433),
    @("F23", 67.08482137424372),
    @("G23", 3.533873050689236),
    @("M23", 34.95462418596033),
    @("N23", 16.85444230536397),
    @("B24", 8.044353639435879),
    @("D24", 11.04839420432523),
    @("E24", 14.27048875209654),
    @("F24", 63.02587081902524),
    @("G24", 3.554326200066487),
    @("M24", 32.84103511870809),
    @("N24", 16.92873582065655),
    @("B25", 7.886967587280284),
    @("D25", 10.35539527624614),
    @("E25", 13.19653530915206),
    @("F25", 58.50506007514766),
    @("G25", 3.577479272937054),
    @("M25", 30.45727794898439),
    @("N25", 17.01799266548944)
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
